$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for price cells whose values would otherwise be
# auto-parsed as numbers (losing exact decimal text / precision).
$textFormatCells = @(
    "D5",
    "D8",
    "D9",
    "D10",
    "D11",
    "D13",
    "D15",
    "D16",
    "D18",
    "D21",
    "D22",
    "D24",
    "D25",
    "D26",
    "D27",
    "D28",
    "D30",
    "D35",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D43",
    "D44",
    "D46",
    "D47",
    "D48"
)
foreach ($cellRef in $textFormatCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply updated coin / link / price / volume values.
$cellValues = @{
    "D2" = '35.595.75'
    "E2" = '  +1.78%  '
    "D3" = '1.908.26'
    "E3" = '  +3.59%  '
    "E4" = '  +0.87%  '
    "D5" = '246.57'
    "E5" = '  +5.95%  '
    "E6" = '  +2.21%  '
    "E7" = '  +0.73%  '
    "D8" = '42.38'
    "E8" = '  +2.45%  '
    "D9" = '0.339'
    "E9" = '  +3.40%  '
    "D10" = '0.0706'
    "E10" = '  +1.83%  '
    "D11" = '0.0997'
    "E11" = '  +1.79%  '
    "D12" = '2.184.95'
    "E12" = '  +3.57%  '
    "D13" = '12.45'
    "E13" = '  +8.86%  '
    "D14" = '1.926.72'
    "E14" = '  +3.78%  '
    "D15" = '0.691'
    "E15" = '  +2.96%  '
    "D16" = '4.86'
    "E16" = '  +3.76%  '
    "D17" = '35.570.19'
    "E17" = '  +1.54%  '
    "D18" = '72.02'
    "E18" = '  +2.97%  '
    "E19" = '  +2.64%  '
    "E20" = '  +1.64%  '
    "D21" = '12.45'
    "E21" = '  +2.41%  '
    "D22" = '4.93'
    "E22" = '  +3.43%  '
    "E23" = '  +0.87%  '
    "D24" = '2.29'
    "E24" = '  -1.05%  '
    "D25" = '171.98'
    "E25" = '  +0.44%  '
    "D26" = '2.19'
    "E26" = '  +26.66%  '
    "D27" = '8.56'
    "E27" = '  +9.06%  '
    "D28" = '17.99'
    "E28" = '  +3.08%  '
    "E29" = '  +1.23%  '
    "D30" = '0.982'
    "E30" = '  +30.10%  '
    "E31" = '  +3.78%  '
    "E32" = '  +2.14%  '
    "E33" = '  +0.72%  '
    "D35" = '1.73'
    "E35" = '  +5.96%  '
    "E36" = '  +3.37%  '
    "D37" = '1.31'
    "E37" = '  +5.00%  '
    "B38" = 'MultiversX'
    "C38" = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    "D38" = '54.33'
    "E38" = '  +58.99%  '
    "B39" = 'ARBITRUM'
    "C39" = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    "D39" = '1.11'
    "E39" = '  +4.94%  '
    "B40" = 'VeChain'
    "C40" = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    "D40" = '0.0206'
    "E40" = '  +3.51%  '
    "B41" = 'Aave'
    "C41" = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
    "D41" = '91.57'
    "E41" = '  +1.74%  '
    "D42" = '1.357.56'
    "E42" = '  +1.13%  '
    "D43" = '15.37'
    "E43" = '  +6.03%  '
    "D44" = '0.0593'
    "E45" = '  +3.87%  '
    "B46" = 'HuobiToken'
    "C46" = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    "D46" = '2.44'
    "E46" = '  +1.31%  '
    "B47" = 'Gas'
    "C47" = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
    "D47" = '12.64'
    "E47" = '  +11.02%  '
    "D48" = '2.77'
    "E48" = '  +0.67%  '
    "E49" = '  +5.67%  '
    "D50" = '2.093.78'
    "E50" = '  +3.51%  '
    "E51" = '  +2.61%  '
}
foreach ($cellRef in $cellValues.Keys) {
    $ws.Range($cellRef).Value = $cellValues[$cellRef]
}
